$wb = $excel.ActiveWorkbook

# Rename the first two worksheets
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "FlowDict_Good"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "StockDict_Good"

# FlowDict_Good (sheet1): update column I (values slightly change due to recomputation
# with different floating-point rounding noise picked up on re-save)
$ws1.Cells.Item(2, 9).Value = 4.999999999999996
$ws1.Cells.Item(3, 9).Value = 5.499999999999996
$ws1.Cells.Item(4, 9).Value = 5.999999999999996
$ws1.Cells.Item(5, 9).Value = 6.499999999999995
$ws1.Cells.Item(6, 9).Value = 6.999999999999995
$ws1.Cells.Item(7, 9).Value = 7.499999999999994
$ws1.Cells.Item(8, 9).Value = 7.999999999999994
$ws1.Cells.Item(9, 9).Value = 8.499999999999993
$ws1.Cells.Item(10, 9).Value = 8.999999999999993
$ws1.Cells.Item(11, 9).Value = 9.499999999999993
$ws1.Cells.Item(12, 9).Value = 9.999999999999993
$ws1.Cells.Item(13, 9).Value = 0.4999999999999996
$ws1.Cells.Item(14, 9).Value = 10.99999999999999
$ws1.Cells.Item(15, 9).Value = 11.49999999999999
$ws1.Cells.Item(16, 9).Value = 11.99999999999999
$ws1.Cells.Item(17, 9).Value = 12.49999999999999
$ws1.Cells.Item(18, 9).Value = 12.99999999999999
$ws1.Cells.Item(19, 9).Value = 13.49999999999999
$ws1.Cells.Item(20, 9).Value = 13.99999999999999
$ws1.Cells.Item(21, 9).Value = 14.49999999999999
$ws1.Cells.Item(22, 9).Value = 14.99999999999999

# StockDict_Good (sheet2): update columns D and E with the same kind of tiny
# floating-point noise adjustments
$ws2.Cells.Item(2, 4).Value = 4.999999999999996
$ws2.Cells.Item(2, 5).Value = 4.999999999999996
$ws2.Cells.Item(3, 4).Value = 5.499999999999996
$ws2.Cells.Item(3, 5).Value = 10.49999999999999
$ws2.Cells.Item(4, 4).Value = 5.999999999999996
$ws2.Cells.Item(4, 5).Value = 16.49999999999999
$ws2.Cells.Item(5, 4).Value = 6.499999999999995
$ws2.Cells.Item(5, 5).Value = 22.99999999999999
$ws2.Cells.Item(6, 4).Value = 6.999999999999995
$ws2.Cells.Item(6, 5).Value = 29.99999999999998
$ws2.Cells.Item(7, 4).Value = 7.499999999999994
$ws2.Cells.Item(7, 5).Value = 37.49999999999997
$ws2.Cells.Item(8, 4).Value = 7.999999999999994
$ws2.Cells.Item(8, 5).Value = 45.49999999999996
$ws2.Cells.Item(9, 4).Value = 8.499999999999993
$ws2.Cells.Item(9, 5).Value = 53.99999999999996
$ws2.Cells.Item(10, 4).Value = 8.999999999999993
$ws2.Cells.Item(10, 5).Value = 62.99999999999995
$ws2.Cells.Item(11, 4).Value = 9.499999999999993
$ws2.Cells.Item(11, 5).Value = 72.49999999999994
$ws2.Cells.Item(12, 4).Value = 9.999999999999993
$ws2.Cells.Item(12, 5).Value = 82.49999999999994
$ws2.Cells.Item(13, 4).Value = 0.4999999999999996
$ws2.Cells.Item(13, 5).Value = 82.99999999999994
$ws2.Cells.Item(14, 4).Value = 10.99999999999999
$ws2.Cells.Item(14, 5).Value = 93.99999999999993
$ws2.Cells.Item(15, 4).Value = 11.49999999999999
$ws2.Cells.Item(15, 5).Value = 105.4999999999999
$ws2.Cells.Item(16, 4).Value = 11.99999999999999
$ws2.Cells.Item(16, 5).Value = 117.4999999999999
$ws2.Cells.Item(17, 4).Value = 12.49999999999999
$ws2.Cells.Item(17, 5).Value = 129.9999999999999
$ws2.Cells.Item(18, 4).Value = 12.99999999999999
$ws2.Cells.Item(18, 5).Value = 142.9999999999999
$ws2.Cells.Item(19, 4).Value = 13.49999999999999
$ws2.Cells.Item(19, 5).Value = 156.4999999999999
$ws2.Cells.Item(20, 4).Value = 13.99999999999999
$ws2.Cells.Item(20, 5).Value = 170.4999999999999
$ws2.Cells.Item(21, 4).Value = 14.49999999999999
$ws2.Cells.Item(21, 5).Value = 184.9999999999999
$ws2.Cells.Item(22, 4).Value = 14.99999999999999
$ws2.Cells.Item(22, 5).Value = 199.9999999999999
